# Scheduled market-data refresh: update price/profit columns (H-N) on each
# Leve-profit worksheet to the latest Universalis averages.
$wb = $excel.ActiveWorkbook

# --- ALC sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 24900
$ws.Range("I7").Value = 24900
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 24900
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -24788
$ws.Range("N7").ClearContents()
$ws.Range("H9").Value = 29
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 29
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 29
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -367
$ws.Range("H14").Value = 24900
$ws.Range("I14").Value = 24900
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 24900
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -24709
$ws.Range("N14").ClearContents()
$ws.Range("H28").Value = 1232.4
$ws.Range("I28").Value = 1122.5555
$ws.Range("K28").Value = 1122.5555
$ws.Range("M28").Value = -637.5554999999999
$ws.Range("H40").Value = 33518.39
$ws.Range("J40").Value = 31301.938
$ws.Range("L40").Value = 31301.938
$ws.Range("N40").Value = -31651.938
$ws.Range("H46").Value = 7600
$ws.Range("J46").Value = 7600
$ws.Range("L46").Value = 22800
$ws.Range("N46").Value = -23038
$ws.Range("H51").Value = 8781.956
$ws.Range("I51").Value = 4373.75
$ws.Range("K51").Value = 4373.75
$ws.Range("M51").Value = -3889.75
$ws.Range("H60").Value = 7600
$ws.Range("J60").Value = 7600
$ws.Range("L60").Value = 22800
$ws.Range("N60").Value = -23768
$ws.Range("H86").Value = 83335670
$ws.Range("I86").Value = 142859200
$ws.Range("J86").Value = 2718.8
$ws.Range("K86").Value = 142859200
$ws.Range("L86").Value = 2718.8
$ws.Range("M86").Value = -142858077
$ws.Range("N86").Value = -4964.8
$ws.Range("H87").Value = 72942.86
$ws.Range("J87").Value = 72942.86
$ws.Range("L87").Value = 72942.86
$ws.Range("N87").Value = -75438.86
$ws.Range("H89").Value = 83335670
$ws.Range("I89").Value = 142859200
$ws.Range("J89").Value = 2718.8
$ws.Range("K89").Value = 714296000
$ws.Range("L89").Value = 13594
$ws.Range("M89").Value = -714290384
$ws.Range("N89").Value = -24826
$ws.Range("H90").Value = 72942.86
$ws.Range("J90").Value = 72942.86
$ws.Range("L90").Value = 218828.58
$ws.Range("N90").Value = -231308.58
$ws.Range("H92").Value = 50000220
$ws.Range("I92").Value = 62500132
$ws.Range("J92").Value = 577
$ws.Range("K92").Value = 62500132
$ws.Range("L92").Value = 577
$ws.Range("M92").Value = -62498884
$ws.Range("N92").Value = -3073
$ws.Range("H96").Value = 660.9
$ws.Range("I96").Value = 660.9
$ws.Range("K96").Value = 1982.7
$ws.Range("M96").Value = -609.6999999999998
$ws.Range("H98").Value = 1595.25
$ws.Range("I98").Value = 1283.1052
$ws.Range("K98").Value = 1283.1052
$ws.Range("M98").Value = 214.8948
$ws.Range("H107").Value = 1118.8
$ws.Range("I107").Value = 1132.1111
$ws.Range("K107").Value = 1132.1111
$ws.Range("M107").Value = 787.8888999999999
$ws.Range("H122").Value = 1595.25
$ws.Range("I122").Value = 1283.1052
$ws.Range("K122").Value = 3849.3156
$ws.Range("M122").Value = -1399.3156
$ws.Range("H125").Value = 2346.1538
$ws.Range("J125").Value = 3750
$ws.Range("L125").Value = 33750
$ws.Range("N125").Value = -38670
$ws.Range("H137").Value = 2520.8823
$ws.Range("J137").Value = 3599.6667
$ws.Range("L137").Value = 10799.0001
$ws.Range("N137").Value = -15899.0001

# --- ARM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30100.445
$ws.Range("I32").Value = 30181.23
$ws.Range("K32").Value = 30181.23
$ws.Range("M32").Value = -29894.23
$ws.Range("H88").Value = 33531.312
$ws.Range("I88").Value = 640.6
$ws.Range("J88").Value = 48481.637
$ws.Range("K88").Value = 640.6
$ws.Range("L88").Value = 48481.637
$ws.Range("M88").Value = -234.6
$ws.Range("N88").Value = -49293.637
$ws.Range("H91").Value = 33531.312
$ws.Range("I91").Value = 640.6
$ws.Range("J91").Value = 48481.637
$ws.Range("K91").Value = 640.6
$ws.Range("L91").Value = 48481.637
$ws.Range("M91").Value = 763.4
$ws.Range("N91").Value = -51289.637
$ws.Range("H110").Value = 928749.0600000001
$ws.Range("I110").Value = 1075072.9
$ws.Range("K110").Value = 1075072.9
$ws.Range("M110").Value = -1073027.9

# --- BSM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 15985.6
$ws.Range("I97").Value = 14982
$ws.Range("K97").Value = 14982
$ws.Range("M97").Value = -13991
$ws.Range("H105").Value = 3812.6
$ws.Range("I105").Value = 3817.7036
$ws.Range("J105").Value = 3766.6667
$ws.Range("K105").Value = 3817.7036
$ws.Range("L105").Value = 3766.6667
$ws.Range("M105").Value = -2070.7036
$ws.Range("N105").Value = -7260.6667

# --- CRP sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 828
$ws.Range("I16").Value = 641.3570999999999
$ws.Range("J16").Value = 1263.5
$ws.Range("K16").Value = 641.3570999999999
$ws.Range("L16").Value = 1263.5
$ws.Range("M16").Value = -354.3570999999999
$ws.Range("N16").Value = -1837.5
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H99").Value = 9199.4
$ws.Range("I99").Value = 6665.6665
$ws.Range("K99").Value = 6665.6665
$ws.Range("M99").Value = -5167.6665
$ws.Range("H113").Value = 828
$ws.Range("I113").Value = 641.3570999999999
$ws.Range("J113").Value = 1263.5
$ws.Range("K113").Value = 641.3570999999999
$ws.Range("L113").Value = 1263.5
$ws.Range("M113").Value = 1528.6429
$ws.Range("N113").Value = -5603.5
$ws.Range("H126").Value = 9199.4
$ws.Range("I126").Value = 6665.6665
$ws.Range("K126").Value = 19996.9995
$ws.Range("M126").Value = -17526.9995

# --- CUL sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1420.1034
$ws.Range("I129").Value = 852.5909
$ws.Range("J129").Value = 3203.7144
$ws.Range("K129").Value = 2557.7727
$ws.Range("L129").Value = 9611.143199999999
$ws.Range("M129").Value = 2442.2273
$ws.Range("N129").Value = -19611.1432
$ws.Range("H132").Value = 1152.8889
$ws.Range("I132").Value = 910.95
$ws.Range("J132").Value = 1844.1428
$ws.Range("K132").Value = 8198.550000000001
$ws.Range("L132").Value = 16597.2852
$ws.Range("M132").Value = -5668.550000000001
$ws.Range("N132").Value = -21657.2852

# --- GSM sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 145904.58
$ws.Range("I80").Value = 254433.25
$ws.Range("J80").Value = 1199.6666
$ws.Range("K80").Value = 254433.25
$ws.Range("L80").Value = 1199.6666
$ws.Range("M80").Value = -253435.25
$ws.Range("N80").Value = -3195.6666
$ws.Range("H83").Value = 145904.58
$ws.Range("I83").Value = 254433.25
$ws.Range("J83").Value = 1199.6666
$ws.Range("K83").Value = 1272166.25
$ws.Range("L83").Value = 5998.333000000001
$ws.Range("M83").Value = -1267174.25
$ws.Range("N83").Value = -15982.333
$ws.Range("H97").Value = 640.6923
$ws.Range("I97").Value = 302.72726
$ws.Range("K97").Value = 302.72726
$ws.Range("M97").Value = 193.27274
$ws.Range("H102").Value = 19238436
$ws.Range("I102").Value = 31258466
$ws.Range("J102").Value = 6386.2
$ws.Range("K102").Value = 31258466
$ws.Range("L102").Value = 6386.2
$ws.Range("M102").Value = -31256844
$ws.Range("N102").Value = -9630.200000000001
$ws.Range("H126").Value = 3628.1724
$ws.Range("I126").Value = 2067.3333
$ws.Range("K126").Value = 6201.999899999999
$ws.Range("M126").Value = -3731.999899999999
$ws.Range("H132").Value = 3343.825
$ws.Range("I132").Value = 3019.1614
$ws.Range("J132").Value = 4462.1113
$ws.Range("K132").Value = 9057.484199999999
$ws.Range("L132").Value = 13386.3339
$ws.Range("M132").Value = -6527.484199999999
$ws.Range("N132").Value = -18446.3339

# --- LTW sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3810.975
$ws.Range("J7").Value = 4842.9653
$ws.Range("L7").Value = 4842.9653
$ws.Range("N7").Value = -5066.9653
$ws.Range("H17").Value = 108
$ws.Range("I17").Value = 108
$ws.Range("K17").Value = 108
$ws.Range("M17").Value = 62
$ws.Range("H57").Value = 22498.5
$ws.Range("I57").Value = 19997
$ws.Range("K57").Value = 19997
$ws.Range("M57").Value = -19431
$ws.Range("H126").Value = 3810.975
$ws.Range("J126").Value = 4842.9653
$ws.Range("L126").Value = 14528.8959
$ws.Range("N126").Value = -19468.8959
$ws.Range("H136").Value = 4936.3335
$ws.Range("I136").Value = 3426.3684
$ws.Range("K136").Value = 10279.1052
$ws.Range("M136").Value = -7729.1052

# --- WVR sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17992.4
$ws.Range("J41").Value = 17992.4
$ws.Range("L41").Value = 17992.4
$ws.Range("N41").Value = -18772.4
$ws.Range("H53").Value = 15000
$ws.Range("I53").Value = 15000
$ws.Range("K53").Value = 15000
$ws.Range("M53").Value = -14393
$ws.Range("H58").Value = 19999.5
$ws.Range("I58").Value = 19999.5
$ws.Range("K58").Value = 19999.5
$ws.Range("M58").Value = -19691.5
$ws.Range("H81").Value = 995835.2
$ws.Range("J81").Value = 3626.9167
$ws.Range("L81").Value = 7253.8334
$ws.Range("N81").Value = -9375.8334
$ws.Range("H84").Value = 995835.2
$ws.Range("J84").Value = 3626.9167
$ws.Range("L84").Value = 36269.167
$ws.Range("N84").Value = -46877.167
$ws.Range("H126").Value = 1476.6666
$ws.Range("J126").Value = 1303
$ws.Range("L126").Value = 3909
$ws.Range("N126").Value = -8849
$ws.Range("H132").Value = 1451.3256
$ws.Range("I132").Value = 984.9474
$ws.Range("K132").Value = 2954.8422
$ws.Range("M132").Value = -424.8422
